# ----------------------------------------------------------------------------
# Applies the "gh-pages output generated at 456a3b4" update to 上海-漫展信息.xlsx
#  - bumps "想去人数" (F column) want-to-go counters across all four sheets
#  - updates the JZ Club venue address on sheet 2 (演出), row 5
#  - appends a new event row (NW新界动漫游戏展2.0) to sheet 1 (展览), row 36
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- 展览 (sheet 1): update "想去人数" (F column) ---
$ws1.Range("F3").Value = 780
$ws1.Range("F5").Value = 2270
$ws1.Range("F6").Value = 1351
$ws1.Range("F7").Value = 105
$ws1.Range("F8").Value = 806
$ws1.Range("F9").Value = 134
$ws1.Range("F11").Value = 2963
$ws1.Range("F14").Value = 1087
$ws1.Range("F15").Value = 601
$ws1.Range("F17").Value = 122
$ws1.Range("F19").Value = 1030
$ws1.Range("F20").Value = 1030
$ws1.Range("F21").Value = 131
$ws1.Range("F22").Value = 12
$ws1.Range("F23").Value = 155
$ws1.Range("F25").Value = 191
$ws1.Range("F26").Value = 633
$ws1.Range("F27").Value = 599
$ws1.Range("F28").Value = 821
$ws1.Range("F29").Value = 46
$ws1.Range("F31").Value = 1016
$ws1.Range("F32").Value = 5000
$ws1.Range("F33").Value = 472
$ws1.Range("F34").Value = 220
$ws1.Range("F35").Value = 112

# --- 演出 (sheet 2): update "想去人数" (F column) ---
$ws2.Range("F6").Value = 406
$ws2.Range("F10").Value = 3
$ws2.Range("F18").Value = 1774
$ws2.Range("F22").Value = 40
$ws2.Range("F23").Value = 47
$ws2.Range("F26").Value = 643
$ws2.Range("F31").Value = 55
$ws2.Range("F34").Value = 350
$ws2.Range("F35").Value = 350
$ws2.Range("F41").Value = 747

# --- 本地生活 (sheet 3): update "想去人数" (F column) ---
$ws3.Range("F4").Value = 633
$ws3.Range("F5").Value = 408
$ws3.Range("F6").Value = 392

# --- 全部类型 (sheet 4): update "想去人数" (F column) ---
$ws4.Range("F4").Value = 408
$ws4.Range("F5").Value = 780
$ws4.Range("F9").Value = 2270
$ws4.Range("F10").Value = 1351
$ws4.Range("F11").Value = 105
$ws4.Range("F12").Value = 806
$ws4.Range("F13").Value = 3
$ws4.Range("F16").Value = 2963
$ws4.Range("F19").Value = 1087
$ws4.Range("F20").Value = 601
$ws4.Range("F22").Value = 392
$ws4.Range("F24").Value = 1774
$ws4.Range("F26").Value = 122
$ws4.Range("F27").Value = 1030
$ws4.Range("F28").Value = 1030
$ws4.Range("F29").Value = 131
$ws4.Range("F31").Value = 12
$ws4.Range("F32").Value = 155
$ws4.Range("F33").Value = 191
$ws4.Range("F34").Value = 40
$ws4.Range("F35").Value = 633
$ws4.Range("F36").Value = 599
$ws4.Range("F38").Value = 643
$ws4.Range("F39").Value = 821
$ws4.Range("F40").Value = 46
$ws4.Range("F41").Value = 1016
$ws4.Range("F42").Value = 5000
$ws4.Range("F43").Value = 55
$ws4.Range("F44").Value = 472
$ws4.Range("F46").Value = 350
$ws4.Range("F47").Value = 220

# --- 演出 (sheet 2): JZ Club moved venue -> update address text (D5) ---
$ws2.Range("D5").Value = "衡山路八号水塔广场 JZ Club 爵士上海俱乐部"

# --- 展览 (sheet 1): append new event row 36 ---
$ws1.Range("A36").Value = 35
# Column B holds the date as plain text elsewhere in the sheet (e.g. B2:B35),
# so force text formatting before assignment - otherwise Excel autodetects
# "2024-07-05" as a date and stores a date serial instead of the literal string.
$ws1.Range("B36").NumberFormat = "@"
$ws1.Range("B36").Value = "2024-07-05"
$ws1.Range("C36").Value = "上海·NW新界动漫游戏展2.0"
$ws1.Range("D36").Value = "长寿路街道澳门路168号 月星家居"
$ws1.Range("E36").Value = "2024.07.05 10:00-07.07 16:00"
$ws1.Range("F36").Value = 0
$ws1.Range("G36").Value = 59
$ws1.Range("H36").Value = "https://show.bilibili.com/platform/detail.html?id=83923"
$ws1.Range("I36").Value = "//i1.hdslb.com/bfs/openplatform/202404/nUixyAy21712605861917.jpeg"

# Mirror row 35's cell style onto the new A36 cell (bold/centered/bordered "#" column style)
$ws1.Range("A35").Copy()
$ws1.Range("A36").PasteSpecial(-4122)
$excel.CutCopyMode = $false

